$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row (A1:D1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case Spanish connector words (de/del/el/la/los/las) in municipality/state names ---
$ws.Range('B4').Value = 'Pabellón De Arteaga'
$ws.Range('B5').Value = 'Rincón De Romos'
$ws.Range('B18').Value = 'Bejucal De Ocampo'
$ws.Range('B21').Value = 'Chiapa De Corzo'
$ws.Range('B33').Value = 'Mazapa De Madero'
$ws.Range('B36').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B38').Value = 'San Cristóbal De Las Casas'
$ws.Range('B62').Value = 'San Juan De Sabinas'
$ws.Range('A65').Value = 'Ciudad De México'
$ws.Range('B69').Value = 'Cuajimalpa De Morelos'
$ws.Range('A85').Value = 'Estado De México'
$ws.Range('B85').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B86').Value = 'Almoloya De Juárez'
$ws.Range('B89').Value = 'Atizapán De Zaragoza'
$ws.Range('B97').Value = 'Ecatepec De Morelos'
$ws.Range('B101').Value = 'Ixtapan De La Sal'
$ws.Range('B109').Value = 'Naucalpan De Juárez'
$ws.Range('B114').Value = 'San Felipe Del Progreso'
$ws.Range('B125').Value = 'Tlalnepantla De Baz'
$ws.Range('B128').Value = 'Valle De Bravo'
$ws.Range('B136').Value = 'Apaseo El Alto'
$ws.Range('B137').Value = 'Apaseo El Grande'
$ws.Range('B143').Value = 'Jaral Del Progreso'
$ws.Range('B149').Value = 'San Diego De La Unión'
$ws.Range('B151').Value = 'San Francisco Del Rincón'
$ws.Range('B153').Value = 'Valle De Santiago'
$ws.Range('B157').Value = 'Acapulco De Juárez'
$ws.Range('B159').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B162').Value = 'Atlamajalcingo Del Monte'
$ws.Range('B164').Value = 'Ayutla De Los Libres'
$ws.Range('B166').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B167').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B170').Value = 'Coyuca De Benítez'
$ws.Range('B173').Value = 'Cutzamala De Pinzón'
$ws.Range('B177').Value = 'Zihuatanejo De Azueta'
$ws.Range('B179').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B190').Value = 'Técpan De Galeana'
$ws.Range('B191').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B194').Value = 'Tlalixtaquilla De Maldonado'
$ws.Range('B195').Value = 'Tlapa De Comonfort'
$ws.Range('B205').Value = 'Atotonilco El Grande'
$ws.Range('B210').Value = 'Cuautepec De Hinojosa'
$ws.Range('B212').Value = 'Huejutla De Reyes'
$ws.Range('B216').Value = 'Molango De Escamilla'
$ws.Range('B218').Value = 'Pachuca De Soto'
$ws.Range('B219').Value = 'Progreso De Obregón'
$ws.Range('B223').Value = 'Tepehuacán De Guerrero'
$ws.Range('B224').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B229').Value = 'Tulancingo De Bravo'
$ws.Range('B231').Value = 'Zacualtipán De Ángeles'
$ws.Range('B234').Value = 'Encarnación De Díaz'
$ws.Range('B236').Value = 'Lagos De Moreno'
$ws.Range('B237').Value = 'San Cristóbal De La Barranca'
$ws.Range('B238').Value = 'San Martín De Bolaños'
$ws.Range('B285').Value = 'Zacualpan De Amilpas'
$ws.Range('B293').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B295').Value = 'Fresnillo De Trujano'
$ws.Range('B296').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B297').Value = 'Ixtlán De Juárez'
$ws.Range('B298').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B300').Value = 'Mazatlán Villa De Flores'
$ws.Range('B301').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B302').Value = 'Mixistlán De La Reforma'
$ws.Range('B303').Value = 'Oaxaca De Juárez'
$ws.Range('B304').Value = 'Putla Villa De Guerrero'
$ws.Range('B305').Value = 'Reforma De Pineda'
$ws.Range('B309').Value = 'San Antonio De La Cal'
$ws.Range('B353').Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range('B354').Value = 'Tlacolula De Matamoros'
$ws.Range('B355').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B356').Value = 'Villa De Zaachila'
$ws.Range('B373').Value = 'Cuayuca De Andrade'
$ws.Range('B382').Value = 'Huehuetlán El Chico'
$ws.Range('B387').Value = 'Izúcar De Matamoros'
$ws.Range('B405').Value = 'Tepanco De López'
$ws.Range('B410').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B423').Value = 'Cadereyta De Montes'
$ws.Range('B424').Value = 'Jalpan De Serra'
$ws.Range('B425').Value = 'Landa De Matamoros'
$ws.Range('B426').Value = 'Pinal De Amoles'
$ws.Range('B436').Value = 'Ciudad Del Maíz'
$ws.Range('B442').Value = 'San Ciro De Acosta'
$ws.Range('B448').Value = 'Villa De La Paz'
$ws.Range('B481').Value = 'San Pablo Del Monte'
$ws.Range('B494').Value = 'Boca Del Río'
$ws.Range('B502').Value = 'Cosamaloapan De Carpio'
$ws.Range('B508').Value = 'Hueyapan De Ocampo'
$ws.Range('B509').Value = 'Ignacio De La Llave'
$ws.Range('B510').Value = 'Ixhuacán De Los Reyes'
$ws.Range('B511').Value = 'Ixhuatlán De Madero'
$ws.Range('B512').Value = 'Ixhuatlán Del Café'
$ws.Range('B516').Value = 'Juchique De Ferrer'
$ws.Range('B521').Value = 'Martínez De La Torre'
$ws.Range('B525').Value = 'Mixtla De Altamirano'
$ws.Range('B527').Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range('B532').Value = 'Paso Del Macho'
$ws.Range('B561').Value = 'Mezquital Del Oro'
$ws.Range('B563').Value = 'Villa De Cos'

# --- Remove trailing metadata/footnote rows 569:573 ---
$ws.Rows("569:573").Delete()
